$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "258.58"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.03%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.01"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-3.41%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.884"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-8.78%"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.26%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.686"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.36%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8751"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.34%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9634"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "5.11%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1416"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.24%"

$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.03592"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "4.15%"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07178"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.14%"

$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03135"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.42%"

$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09232"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.14%"

$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001550"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.17%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005997"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.65%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.484"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.49%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.222"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-2.44%"

$ws.Range("B19").Value = "One"
$ws.Range("C19").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.01063"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1,646.55%"

$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3145"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.68%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1306"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.29%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.524"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.04%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04217"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.75%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1379"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.09%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001219"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.68%"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-7.26%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001199"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.01%"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001492"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "2.61%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03836"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.39%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.005886"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.75%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1103"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.25%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002199"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.01%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01049"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "5.81%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005490"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "3.78%"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.01%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1090"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "9.11%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002154"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "1.22%"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.01%"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.01%"
